$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 106008.5
$ws.Range("I28").Value = 127011.75
$ws.Range("K28").Value = 127011.75
$ws.Range("M28").Value = -126526.75

$ws.Range("H100").Value = 1075
$ws.Range("I100").Value = 1075
$ws.Range("K100").Value = 1075
$ws.Range("M100").Value = -534

$ws.Range("H112").Value = 2460.95
$ws.Range("J112").Value = 1731.2667
$ws.Range("L112").Value = 5193.800099999999
$ws.Range("N112").Value = -7409.800099999999

$ws.Range("H138").Value = 3019.6292
$ws.Range("I138").Value = 1406.9546
$ws.Range("K138").Value = 4220.8638
$ws.Range("M138").Value = 919.1361999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 4277.154
$ws.Range("J5").Value = 12620.5
$ws.Range("L5").Value = 12620.5
$ws.Range("N5").Value = -12844.5

$ws.Range("H102").Value = 1029901.6
$ws.Range("I102").Value = 1029901.6
$ws.Range("K102").Value = 1029901.6
$ws.Range("M102").Value = -1028279.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 4277.154
$ws.Range("J4").Value = 12620.5
$ws.Range("L4").Value = 12620.5
$ws.Range("N4").Value = -12850.5

$ws.Range("H22").Value = 205.54546
$ws.Range("I22").Value = 212.375
$ws.Range("K22").Value = 212.375
$ws.Range("M22").Value = -39.375

$ws.Range("H105").Value = 3579.2
$ws.Range("I105").Value = 4299.9
$ws.Range("J105").Value = 2137.8
$ws.Range("K105").Value = 4299.9
$ws.Range("L105").Value = 2137.8
$ws.Range("M105").Value = -2552.9
$ws.Range("N105").Value = -5631.8

$ws.Range("H134").Value = 1258.9025
$ws.Range("I134").Value = 1083.5161
$ws.Range("K134").Value = 3250.5483
$ws.Range("M134").Value = -715.5483000000004

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1883.6
$ws.Range("J16").Value = 1769.6
$ws.Range("L16").Value = 1769.6
$ws.Range("N16").Value = -2343.6

$ws.Range("H22").Value = 358.85294
$ws.Range("J22").Value = 347.46155
$ws.Range("L22").Value = 347.46155
$ws.Range("N22").Value = -1047.46155

$ws.Range("H31").Value = 5073.4155
$ws.Range("I31").Value = 22143.5
$ws.Range("J31").Value = 3337.4746
$ws.Range("K31").Value = 22143.5
$ws.Range("L31").Value = 3337.4746
$ws.Range("M31").Value = -21848.5
$ws.Range("N31").Value = -3927.4746

$ws.Range("H34").Value = 5073.4155
$ws.Range("I34").Value = 22143.5
$ws.Range("J34").Value = 3337.4746
$ws.Range("K34").Value = 22143.5
$ws.Range("L34").Value = 3337.4746
$ws.Range("M34").Value = -21941.5
$ws.Range("N34").Value = -3741.4746

$ws.Range("H86").Value = 6329.1
$ws.Range("I86").Value = 7872.5
$ws.Range("J86").Value = 5300.1665
$ws.Range("K86").Value = 7872.5
$ws.Range("L86").Value = 5300.1665
$ws.Range("M86").Value = -6749.5
$ws.Range("N86").Value = -7546.1665

$ws.Range("H89").Value = 6329.1
$ws.Range("I89").Value = 7872.5
$ws.Range("J89").Value = 5300.1665
$ws.Range("K89").Value = 39362.5
$ws.Range("L89").Value = 26500.8325
$ws.Range("M89").Value = -33746.5
$ws.Range("N89").Value = -37732.8325

$ws.Range("H107").Value = 4982.5264
$ws.Range("I107").Value = 3744.6667
$ws.Range("K107").Value = 3744.6667
$ws.Range("M107").Value = -1824.6667

$ws.Range("H113").Value = 1883.6
$ws.Range("J113").Value = 1769.6
$ws.Range("L113").Value = 1769.6
$ws.Range("N113").Value = -6109.6

$ws.Range("H132").Value = 607204.6
$ws.Range("I132").Value = 417282.03
$ws.Range("J132").Value = 1113664.8
$ws.Range("K132").Value = 1251846.09
$ws.Range("L132").Value = 3340994.4
$ws.Range("M132").Value = -1249316.09
$ws.Range("N132").Value = -3346054.4

$ws.Range("H141").Value = 81157.6
$ws.Range("J141").Value = 81157.6
$ws.Range("L141").Value = 81157.6
$ws.Range("N141").Value = -91517.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 14607460
$ws.Range("I4").Value = 26998980
$ws.Range("K4").Value = 80996940
$ws.Range("M4").Value = -80996828

$ws.Range("H55").Value = 799.8
$ws.Range("J55").Value = 1000
$ws.Range("L55").Value = 3000
$ws.Range("N55").Value = -3354

$ws.Range("H86").Value = 581
$ws.Range("I86").Value = 444.66666
$ws.Range("J86").Value = 785.5
$ws.Range("K86").Value = 1333.99998
$ws.Range("L86").Value = 2356.5
$ws.Range("M86").Value = -147.9999800000001
$ws.Range("N86").Value = -4728.5

$ws.Range("H89").Value = 581
$ws.Range("I89").Value = 444.66666
$ws.Range("J89").Value = 785.5
$ws.Range("K89").Value = 4001.99994
$ws.Range("L89").Value = 7069.5
$ws.Range("M89").Value = 1926.00006
$ws.Range("N89").Value = -18925.5

$ws.Range("H126").Value = 666666.3
$ws.Range("I126").Value = 666666.3
$ws.Range("K126").Value = 1999998.9
$ws.Range("M126").Value = -1995058.9

$ws.Range("H129").Value = 2107.1177
$ws.Range("J129").Value = 2308.4167
$ws.Range("L129").Value = 6925.250100000001
$ws.Range("N129").Value = -16925.2501

$ws.Range("H131").Value = 4742.85
$ws.Range("J131").Value = 6002.357
$ws.Range("L131").Value = 18007.071
$ws.Range("N131").Value = -28087.071

$ws.Range("H134").Value = 2159.2222
$ws.Range("I134").Value = 2159.2222
$ws.Range("K134").Value = 6477.6666
$ws.Range("M134").Value = -1407.6666

$ws.Range("H136").Value = 3798
$ws.Range("I136").Value = 3798
$ws.Range("K136").Value = 11394
$ws.Range("M136").Value = -6294

$ws.Range("H138").Value = 3814252
$ws.Range("I138").Value = 5001711
$ws.Range("K138").Value = 15005133
$ws.Range("M138").Value = -14999993

$ws.Range("H139").Value = 2091604.5
$ws.Range("I139").Value = 2230878.2
$ws.Range("K139").Value = 6692634.600000001
$ws.Range("M139").Value = -6687494.600000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("N2").Value = ""

$ws.Range("H22").Value = 627.75
$ws.Range("I22").Value = 436.625
$ws.Range("K22").Value = 436.625
$ws.Range("M22").Value = -141.625

$ws.Range("H27").Value = 627.75
$ws.Range("I27").Value = 436.625
$ws.Range("K27").Value = 436.625
$ws.Range("M27").Value = -329.625

$ws.Range("H95").Value = 30710.6
$ws.Range("J95").Value = 30710.6
$ws.Range("L95").Value = 30710.6
$ws.Range("N95").Value = -36202.6

$ws.Range("H110").Value = 52399.8
$ws.Range("J110").Value = 52399.8
$ws.Range("L110").Value = 52399.8
$ws.Range("N110").Value = -60579.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 8688.6
$ws.Range("J14").Value = 8126.5713
$ws.Range("L14").Value = 8126.5713
$ws.Range("N14").Value = -8462.5713

$ws.Range("H21").Value = 24966.666
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 24966.666
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 24966.666
$ws.Range("M21").Value = ""
$ws.Range("N21").Value = -25436.666

$ws.Range("H35").Value = 24966.666
$ws.Range("I35").Value = 0
$ws.Range("J35").Value = 24966.666
$ws.Range("K35").Value = 0
$ws.Range("L35").Value = 24966.666
$ws.Range("M35").Value = ""
$ws.Range("N35").Value = -25546.666

$ws.Range("H40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("N40").Value = ""

$ws.Range("H42").Value = 0
$ws.Range("J42").Value = 0
$ws.Range("L42").Value = 0
$ws.Range("N42").Value = ""

$ws.Range("H64").Value = 147777
$ws.Range("I64").Value = 147777
$ws.Range("K64").Value = 147777
$ws.Range("M64").Value = -147529

$ws.Range("H67").Value = 147777
$ws.Range("I67").Value = 147777
$ws.Range("K67").Value = 147777
$ws.Range("M67").Value = -146919

$ws.Range("H136").Value = 2649.173
$ws.Range("I136").Value = 2010.3636
$ws.Range("K136").Value = 6031.0908
$ws.Range("M136").Value = -3481.0908
